$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values for rows 30-32 (revised figures for Jun/Jul/already-reported months)
$ws.Cells.Item(30, 2).Value = 140788
$ws.Cells.Item(30, 10).Value = 18883

$ws.Cells.Item(31, 2).Value = 146076
$ws.Cells.Item(31, 10).Value = 19198

$ws.Cells.Item(32, 2).Value = 153029
$ws.Cells.Item(32, 9).Value = 29723
$ws.Cells.Item(32, 10).Value = 20162

# Add new row 33 with the new month (01-08-2021) data.
# Writing the date-like label directly via .Value would be auto-parsed as a
# real date (and would bloat styles.xml with a new number format); instead
# stage the literal text via a formula in a scratch cell, then copy/paste
# it as a value so it lands as a plain shared-string cell, exactly like the
# other "Serie" labels in column A.
$ws.Cells.Item(50, 1).Formula = "=""01-08-2021"""
$ws.Cells.Item(50, 1).Copy()
$ws.Cells.Item(33, 1).PasteSpecial(-4163)
$ws.Cells.Item(50, 1).ClearContents()

$ws.Cells.Item(33, 2).Value = 165196
$ws.Cells.Item(33, 3).Value = 30336
$ws.Cells.Item(33, 4).Value = 17945
$ws.Cells.Item(33, 5).Value = 10984
$ws.Cells.Item(33, 6).Value = 11726
$ws.Cells.Item(33, 7).Value = 13799
$ws.Cells.Item(33, 8).Value = 29677
$ws.Cells.Item(33, 9).Value = 29928
$ws.Cells.Item(33, 10).Value = 20801
